$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header row, pushing all existing
# data rows (old 2..38) down to (new 5..41). This also grows the used
# range / dimension automatically.
$ws.Rows("2:4").Insert()

# Insert() copies formatting down from the header row (bold, centered).
# Reset the new rows back to the plain (unstyled) look used by the rest
# of the data rows before re-applying the date format on column D below.
$ws.Range("A2:T4").ClearFormats()

# --- New row 2: Especial, Región de O'Higgins, 2023-06-15 ---
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Femacal de La Calera"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = (Get-Date -Year 2023 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107001
$ws.Range("J2").Value = "Caqui"
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("Q2").Value = "`$/bandeja 10 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1300
$ws.Range("T2").Value = 10

# --- New row 3: Primera, Región de O'Higgins, 2023-06-15 ---
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = (Get-Date -Year 2023 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107001
$ws.Range("J3").Value = "Caqui"
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1200
$ws.Range("T3").Value = 10

# --- New row 4: Segunda, Región de O'Higgins, 2023-06-15 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D4").Value = (Get-Date -Year 2023 -Month 6 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107001
$ws.Range("J4").Value = "Caqui"
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = "`$/bandeja 10 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 10
